# Update cfb_weather.xlsx with Timestamp 2024-11-04T05:15:57.902315
$wb = $excel.ActiveWorkbook

$wsFBS = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

# Update wind direction (wind_dir_fg) values on the "FBS" sheet
$wsFBS.Range("Q10").Value = "SSW"
$wsFBS.Range("Q11").Value = "N"
$wsFBS.Range("Q14").Value = "NW"
$wsFBS.Range("Q16").Value = "NW"
$wsFBS.Range("Q19").Value = "SSW"
$wsFBS.Range("Q35").Value = "ENE"
$wsFBS.Range("Q36").Value = "ENE"
$wsFBS.Range("Q37").Value = "SSW"
$wsFBS.Range("Q39").Value = "ESE"
$wsFBS.Range("Q40").Value = "ESE"
$wsFBS.Range("Q48").Value = "W"

# Update wind direction values on the "Other" sheet
$wsOther.Range("S32").Value = "NW"
$wsOther.Range("S39").Value = "W"
$wsOther.Range("S40").Value = "W"
$wsOther.Range("S46").Value = "ENE"

# Update the Timestamp column (shared across all rows) on the "FBS" sheet
$wsFBS.Range("AK2:AK49").Value = "2024-11-04T05:15:57.902315"
